# Hide the final (visible) slide of the deck.
# Slide 16 is the last slide that isn't already marked hidden
# (slide 17 — the "adapted from" credits slide — is already hidden).
# Setting SlideShowTransition.Hidden = True serializes as show="0"
# on the <p:sld> element, matching PowerPoint's "Hide Slide" command.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(16)
$s.SlideShowTransition.Hidden = $true
